$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# --- Swap the contents of B2 and C2 -----------------------------------
# Before: B2 = "astegic@123" (hyperlink to mailto:astegic@123)
#         C2 = "654321"      (quoted/text number)
# After : B2 = "654321"      (quoted/text number)
#         C2 = "astegic@123" (hyperlink to mailto:astegic@123)
# Route the move through a scratch cell so the hyperlinked cell keeps its
# identity (and therefore its formatting) as it relocates from B2 to C2.
$ws.Range("B2").Cut($ws.Range("Z1")) | Out-Null
$ws.Range("C2").Cut($ws.Range("B2")) | Out-Null
$ws.Range("Z1").Cut($ws.Range("C2")) | Out-Null
$ws.Range("Z1").Clear() | Out-Null

# --- Re-point the hyperlink from B2 to C2 ------------------------------
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:vbansal@astegic.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:astegic@123") | Out-Null

# Re-apply the canonical Hyperlink style (Add() above can create a
# duplicate style record, so make sure both hyperlinked cells end up
# using the worksheet's normal Hyperlink look).
$ws.Range("A2").Style = "Hyperlink"
$ws.Range("C2").Style = "Hyperlink"

# --- Update the selected cell shown in the sheet view ------------------
$ws.Range("C2").Select() | Out-Null

$wb.Save()
